# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table with refreshed snapshot figures (GitHub Actions data refresh).
#
# D-column values are plain decimal-looking strings (e.g. "240.60") that
# must stay text, exactly as authored in the source file, not become
# floating point numbers. A leading apostrophe forces Excel to store the
# value as text (quote-prefixed); Style = 'Normal' afterwards clears the
# resulting quote-prefix style flag so the cell keeps no explicit style,
# matching how these cells were originally formatted.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''36.060.75'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.45%  '
$ws.Range('D3').Value = '''1.918.01'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -4.48%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''240.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.77%  '
$ws.Range('D6').Value = '''0.602'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.76%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '''55.49'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -11.34%  '
$ws.Range('D9').Value = '''0.361'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.31%  '
$ws.Range('D10').Value = '''54.69'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.03%  '
$ws.Range('D11').Value = '''0.0822'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.42%  '
$ws.Range('D12').Value = '''0.103'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.87%  '
$ws.Range('D13').Value = '''2.196.66'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.58%  '
$ws.Range('D14').Value = '''0.803'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -9.44%  '
$ws.Range('D15').Value = '''20.73'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -8.35%  '
$ws.Range('D16').Value = '''13.13'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -7.05%  '
$ws.Range('D17').Value = '''5.15'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.52%  '
$ws.Range('D18').Value = '''1.915.25'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.60%  '
$ws.Range('D19').Value = '''35.977.26'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.46%  '
$ws.Range('D20').Value = '''68.40'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.95%  '
$ws.Range('D21').Value = '''0.0₃0850'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.50%  '
$ws.Range('D22').Value = '''225.45'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.64%  '
$ws.Range('D23').Value = '''4.92'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.56%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').Value = '''2.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.72%  '
$ws.Range('E26').Value = '  -3.12%  '
$ws.Range('D27').Value = '''9.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.81%  '
$ws.Range('D28').Value = '''161.49'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.08%  '
$ws.Range('D29').Value = '''19.07'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.52%  '
$ws.Range('D30').Value = '''0.119'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -12.41%  '
$ws.Range('E31').Value = '  -4.02%  '
$ws.Range('E32').Value = '  -5.35%  '
$ws.Range('E33').Value = '  -8.54%  '
$ws.Range('D34').Value = '''0.0620'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.84%  '
$ws.Range('D35').Value = '''4.24'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.09%  '
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('D38').Value = '''5.89'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -9.57%  '
$ws.Range('D39').Value = '''2.11'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -10.71%  '
$ws.Range('D40').Value = '''2.86'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -11.04%  '
$ws.Range('D41').Value = '''0.0952'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.44%  '
$ws.Range('E42').Value = '  -2.93%  '
$ws.Range('E43').Value = '  -9.74%  '
$ws.Range('E44').Value = '  -4.67%  '
$ws.Range('D45').Value = '''15.45'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.96%  '
$ws.Range('D46').Value = '''1.328.33'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.38%  '
$ws.Range('D47').Value = '''1.02'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -9.41%  '
$ws.Range('D48').Value = '''86.48'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -7.62%  '
$ws.Range('D49').Value = '''7.11'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.34%  '
$ws.Range('D50').Value = '''2.80'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.14%  '
$ws.Range('D51').Value = '''44.57'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.34%  '
